$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly roll-forward: rows 578-703 get new / shifted Fecha (D), Volumen (J),
# Precio minimo (K), Precio maximo (L), Precio promedio ponderado (M) and
# Precio $/Kg (P) values. Rows 578-579 carry the newly observed week; every
# other row (580-703) inherits the values the row two places above it held
# before this edit (the historical series shifts down by one week / two
# quality rows).
$values = @{
  578 = @(44711, 2540, 800, 900, 850, 850)
  579 = @(44711, 1360, 700, 750, 725, 725)
  580 = @(44579, 2400, 650, 700, 675, 675)
  581 = @(44579, 1360, 550, 600, 575, 575)
  582 = @(44342, 3360, 650, 700, 675, 675)
  583 = @(44342, 1800, 500, 550, 525, 525)
  584 = @(44160, 2800, 450, 500, 475, 475)
  585 = @(44160, 1600, 350, 400, 375, 375)
  586 = @(44243, 2400, 650, 700, 675, 675)
  587 = @(44243, 1440, 550, 600, 575, 575)
  588 = @(44251, 3400, 900, 1000, 950, 950)
  589 = @(44251, 2000, 750, 800, 775, 775)
  590 = @(44279, 3400, 700, 750, 725, 725)
  591 = @(44279, 1800, 600, 650, 625, 625)
  592 = @(44610, 2500, 850, 900, 875, 875)
  593 = @(44610, 1300, 750, 800, 775, 775)
  594 = @(44634, 2400, 950, 1000, 975, 975)
  595 = @(44634, 1260, 850, 900, 875, 875)
  596 = @(44519, 3360, 550, 600, 575, 575)
  597 = @(44519, 1600, 450, 500, 475, 475)
  598 = @(44641, 2500, 950, 1000, 975, 975)
  599 = @(44641, 1240, 850, 900, 875, 875)
  600 = @(44505, 3360, 600, 700, 650, 650)
  601 = @(44505, 1600, 500, 550, 525, 525)
  602 = @(44372, 3360, 600, 700, 650, 650)
  603 = @(44372, 1600, 500, 550, 525, 525)
  604 = @(44267, 3000, 850, 900, 875, 875)
  605 = @(44267, 1400, 750, 800, 775, 775)
  606 = @(44669, 2400, 800, 900, 850, 850)
  607 = @(44669, 1200, 700, 750, 725, 725)
  608 = @(44477, 3400, 650, 700, 675, 675)
  609 = @(44477, 1520, 500, 600, 550, 550)
  610 = @(44671, 2600, 850, 900, 875, 875)
  611 = @(44671, 1400, 750, 800, 775, 775)
  612 = @(44474, 2200, 650, 700, 675, 675)
  613 = @(44474, 1340, 550, 600, 575, 575)
  614 = @(44571, 3000, 600, 700, 650, 650)
  615 = @(44571, 1480, 500, 550, 525, 525)
  616 = @(44631, 2460, 950, 1000, 975, 975)
  617 = @(44631, 1300, 850, 900, 875, 875)
  618 = @(44672, 2000, 850, 900, 875, 875)
  619 = @(44672, 1300, 750, 800, 775, 775)
  620 = @(44490, 2000, 600, 700, 650, 650)
  621 = @(44490, 1300, 500, 550, 525, 525)
  622 = @(44679, 2000, 850, 900, 875, 875)
  623 = @(44679, 1300, 750, 800, 775, 775)
  624 = @(44369, 2000, 600, 700, 650, 650)
  625 = @(44369, 1400, 500, 550, 525, 525)
  626 = @(44365, 3400, 600, 700, 650, 650)
  627 = @(44365, 1600, 500, 550, 525, 525)
  628 = @(44603, 2460, 750, 800, 775, 775)
  629 = @(44603, 1280, 650, 700, 675, 675)
  630 = @(44575, 3320, 650, 700, 675, 675)
  631 = @(44575, 1600, 550, 600, 575, 575)
  632 = @(44427, 2200, 650, 700, 675, 675)
  633 = @(44427, 1400, 550, 600, 575, 575)
  634 = @(44587, 2800, 750, 800, 775, 775)
  635 = @(44587, 1500, 650, 700, 675, 675)
  636 = @(44565, 2360, 600, 700, 650, 650)
  637 = @(44565, 1340, 500, 550, 525, 525)
  638 = @(44447, 3000, 600, 700, 650, 650)
  639 = @(44447, 1600, 500, 550, 525, 525)
  640 = @(44445, 3200, 600, 700, 650, 650)
  641 = @(44445, 1500, 500, 550, 525, 525)
  642 = @(44533, 3400, 550, 600, 575, 575)
  643 = @(44533, 1600, 450, 500, 475, 475)
  644 = @(44523, 2400, 600, 700, 650, 650)
  645 = @(44523, 1360, 500, 550, 525, 525)
  646 = @(44601, 2500, 750, 800, 775, 775)
  647 = @(44601, 1360, 650, 700, 675, 675)
  648 = @(44343, 2400, 650, 700, 675, 675)
  649 = @(44343, 1340, 500, 550, 525, 525)
  650 = @(44462, 2000, 600, 700, 650, 650)
  651 = @(44462, 1300, 500, 550, 525, 525)
  652 = @(44159, 2000, 450, 500, 475, 475)
  653 = @(44159, 1400, 350, 400, 375, 375)
  654 = @(44629, 2400, 950, 1000, 975, 975)
  655 = @(44629, 1200, 850, 900, 875, 875)
  656 = @(44536, 2500, 550, 600, 575, 575)
  657 = @(44536, 1500, 450, 500, 475, 475)
  658 = @(44606, 2400, 750, 800, 775, 775)
  659 = @(44606, 1100, 650, 700, 675, 675)
  660 = @(44594, 2400, 750, 800, 775, 775)
  661 = @(44594, 1360, 650, 700, 675, 675)
  662 = @(44377, 3400, 600, 700, 650, 650)
  663 = @(44377, 1800, 500, 550, 525, 525)
  664 = @(44417, 3200, 650, 700, 675, 675)
  665 = @(44417, 1480, 550, 600, 575, 575)
  666 = @(44706, 2640, 850, 900, 875, 875)
  667 = @(44706, 1520, 750, 800, 775, 775)
  668 = @(44664, 2600, 800, 900, 850, 850)
  669 = @(44664, 1340, 700, 750, 725, 725)
  670 = @(44566, 3200, 600, 700, 650, 650)
  671 = @(44566, 1600, 500, 550, 525, 525)
  672 = @(44344, 3400, 650, 700, 675, 675)
  673 = @(44344, 1600, 500, 550, 525, 525)
  674 = @(44351, 3400, 600, 700, 650, 650)
  675 = @(44351, 1600, 500, 550, 525, 525)
  676 = @(44508, 2500, 600, 700, 650, 650)
  677 = @(44508, 1500, 500, 550, 525, 525)
  678 = @(44600, 2200, 700, 800, 750, 750)
  679 = @(44600, 1260, 600, 650, 625, 625)
  680 = @(44323, 3440, 650, 700, 675, 675)
  681 = @(44323, 1660, 500, 550, 525, 525)
  682 = @(44515, 2600, 600, 650, 625, 625)
  683 = @(44515, 1500, 500, 550, 525, 525)
  684 = @(44602, 2100, 750, 800, 775, 775)
  685 = @(44602, 1200, 650, 700, 675, 675)
  686 = @(44326, 3220, 650, 700, 675, 675)
  687 = @(44326, 1460, 500, 550, 525, 525)
  688 = @(44692, 2600, 750, 800, 775, 775)
  689 = @(44692, 1440, 650, 700, 675, 675)
  690 = @(44165, 2800, 450, 500, 475, 475)
  691 = @(44165, 1560, 350, 400, 375, 375)
  692 = @(44655, 2400, 750, 800, 775, 775)
  693 = @(44655, 1280, 650, 700, 675, 675)
  694 = @(44315, 2500, 650, 700, 675, 675)
  695 = @(44315, 1360, 500, 550, 525, 525)
  696 = @(44448, 2000, 600, 700, 650, 650)
  697 = @(44448, 1300, 500, 550, 525, 525)
  698 = @(44263, 3200, 850, 900, 875, 875)
  699 = @(44263, 1600, 750, 800, 775, 775)
  700 = @(44648, 2500, 850, 900, 875, 875)
  701 = @(44648, 1260, 750, 800, 775, 775)
  702 = @(44376, 2460, 600, 700, 650, 650)
  703 = @(44376, 1400, 500, 550, 525, 525)
}

foreach ($row in $values.Keys) {
  $v = $values[$row]
  $ws.Cells.Item($row, 4).Value2  = $v[0]   # D: Fecha
  $ws.Cells.Item($row, 10).Value2 = $v[1]   # J: Volumen
  $ws.Cells.Item($row, 11).Value2 = $v[2]   # K: Precio minimo
  $ws.Cells.Item($row, 12).Value2 = $v[3]   # L: Precio maximo
  $ws.Cells.Item($row, 13).Value2 = $v[4]   # M: Precio promedio ponderado
  $ws.Cells.Item($row, 16).Value2 = $v[5]   # P: Precio $/Kg
}

# Two brand new rows (704-705) are appended, duplicating what used to be the
# oldest pair of rows (702-703) before the roll-forward above overwrote them.
$newRows = @(
  @{ Row = 704; A = 8; B = "Terminal La Palmera de La Serena"; C = "Coquimbo"; D = 44442; E = 4; F = 100112023; G = "Brócoli"; H = "Sin especificar"; I = "Primera"; J = 3200; K = 600; L = 700; M = 650; N = "`$/unidad"; O = "Provincia del Elquí"; P = 650; Q = 1; R = "Hortaliza" },
  @{ Row = 705; A = 8; B = "Terminal La Palmera de La Serena"; C = "Coquimbo"; D = 44442; E = 4; F = 100112023; G = "Brócoli"; H = "Sin especificar"; I = "Segunda"; J = 1600; K = 500; L = 550; M = 525; N = "`$/unidad"; O = "Provincia del Elquí"; P = 525; Q = 1; R = "Hortaliza" }
)

$dateStyleSource = $ws.Range("D703")

foreach ($r in $newRows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value2  = $r.A
  $ws.Cells.Item($row, 2).Value2  = $r.B
  $ws.Cells.Item($row, 3).Value2  = $r.C
  $dCell = $ws.Cells.Item($row, 4)
  $dCell.Value2 = $r.D
  $dCell.NumberFormat = $dateStyleSource.NumberFormat
  $ws.Cells.Item($row, 5).Value2  = $r.E
  $ws.Cells.Item($row, 6).Value2  = $r.F
  $ws.Cells.Item($row, 7).Value2  = $r.G
  $ws.Cells.Item($row, 8).Value2  = $r.H
  $ws.Cells.Item($row, 9).Value2  = $r.I
  $ws.Cells.Item($row, 10).Value2 = $r.J
  $ws.Cells.Item($row, 11).Value2 = $r.K
  $ws.Cells.Item($row, 12).Value2 = $r.L
  $ws.Cells.Item($row, 13).Value2 = $r.M
  $ws.Cells.Item($row, 14).Value2 = $r.N
  $ws.Cells.Item($row, 15).Value2 = $r.O
  $ws.Cells.Item($row, 16).Value2 = $r.P
  $ws.Cells.Item($row, 17).Value2 = $r.Q
  $ws.Cells.Item($row, 18).Value2 = $r.R
}
